$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old _GoBack bookmark (currently sitting in the empty
#    paragraph right before "Usuarios: DNI, contraseña y tipo").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Locate the last paragraph of the document (the "INSERT INTO USUARIOS"
#    paragraph) - everything new gets appended after it, right before
#    the section break.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Grab a template list paragraph (numId=7 numbered item) so the two new
# numbered bullets can continue the very same list instead of starting a
# brand new one.
$listModelPara = $d.Paragraphs.Item(86)
$listTemplate = $listModelPara.Range.ListFormat.ListTemplate

# ---------------------------------------------------------------------------
# A6: new empty "Prrafodelista" paragraph
# ---------------------------------------------------------------------------
$lastPara.Range.InsertParagraphAfter()
$pA6 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pA6.Range.Style = "Prrafodelista"
$pA6.Range.HighlightColorIndex = 0

# ---------------------------------------------------------------------------
# A7: numbered paragraph "Formulario contacto: identificador, nombre, mail,
#     texto" + line break. We place the new _GoBack bookmark at the very end
#     of this paragraph's text, exactly like it used to sit in the old
#     location.
# ---------------------------------------------------------------------------
$pA6.Range.InsertParagraphAfter()
$pA7 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pA7.Range.Style = "Prrafodelista"
$pA7.Range.HighlightColorIndex = 0
$pA7.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
# Trailing "Z" placeholder lets us anchor a collapsed bookmark safely (the
# COM host mishandles a bookmark placed exactly on a paragraph mark); we
# delete the placeholder right after creating the bookmark.
$pA7.Range.Text = "Formulario contacto: identificador, nombre, mail, textoZ"
$bmPos = $pA7.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range($pA7.Range.End - 2, $pA7.Range.End - 1)
$placeholder.Delete()

$brRange = $d.Range($pA7.Range.End - 1, $pA7.Range.End - 1)
$brRange.InsertAfter([char]11)

# ---------------------------------------------------------------------------
# A8: new empty "Prrafodelista" paragraph
# ---------------------------------------------------------------------------
$pA7.Range.InsertParagraphAfter()
$pA8 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pA8.Range.Style = "Prrafodelista"
$pA8.Range.HighlightColorIndex = 0

# ---------------------------------------------------------------------------
# A9: "create table if not exists CONTACTO (...)" paragraph, highlighted
#     yellow.
# ---------------------------------------------------------------------------
$pA8.Range.InsertParagraphAfter()
$pA9 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pA9.Range.Style = "Prrafodelista"
$pA9.Range.Text = "create table if not exists CONTACTO (ID_FORMULARIO INT NOT NULL AUTO_INCREMENT, NOMBRE TEXT, EMAIL TEXT, MENSAJE TEXT, PRIMARY KEY(ID_FORMULARIO));"
$pA9.Range.HighlightColorIndex = 7

Write-Output "Done. Paragraphs=$($d.Paragraphs.Count)"
